# Auto-generated edit script applying the diff changes to 北京-漫展信息.xlsx
# Maps sheet1..sheet4 (OOXML part order) onto Worksheets.Item(1..4):
#   1 = 展览 (Exhibitions), 2 = 演出 (Performances),
#   3 = 本地生活 (Local Life), 4 = 全部类型 (All Types)
$wb = $excel.ActiveWorkbook

# --- 展览 (Worksheets.Item(1)) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range('F3').Value = 11  # was 10
$ws.Range('F5').Value = 1324  # was 1314
$ws.Range('F6').Value = 7649  # was 7643
$ws.Range('F8').Value = 109  # was 110
$ws.Range('F10').Value = 8446  # was 8440
$ws.Range('F14').Value = 5655  # was 5652
$ws.Range('F15').Value = 59  # was 58
$ws.Range('F16').Value = 2616  # was 2608
$ws.Range('F17').Value = 1136  # was 1131
$ws.Range('F18').Value = 4596  # was 4594
$ws.Range('F19').Value = 344  # was 342
$ws.Range('F22').Value = 33  # was 32
$ws.Range('F23').Value = 531  # was 528
$ws.Range('F24').Value = 3507  # was 3487
$ws.Range('F25').Value = 55  # was 44
$ws.Range('F27').Value = 21  # was 20
$ws.Range('F29').Value = 3000  # was 2976
$ws.Range('F30').Value = 31  # was 17
$ws.Range('F31').Value = 99  # was 82
$ws.Range('F32').Value = 343  # was 340
$ws.Range('F34').Value = 310  # was 307
$ws.Range('F35').Value = 351  # was 304
$ws.Range('F38').Value = 882  # was 879
$ws.Range('F39').Value = 1762  # was 1721
$ws.Range('F42').Value = 18  # was 16
$ws.Range('F43').Value = 2845  # was 2749
$ws.Range('F45').Value = 2284  # was 2283
$ws.Range('F47').Value = 30  # was 29

# --- 演出 (Worksheets.Item(2)) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range('F2').Value = 104  # was 103
$ws.Range('F3').Value = 123  # was 121
$ws.Range('F4').Value = 7  # was 6

# --- 本地生活 (Worksheets.Item(3)) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range('F3').Value = 1324  # was 1322

# --- 全部类型 (Worksheets.Item(4)) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range('F3').Value = 1324  # was 1322
$ws.Range('F4').Value = 11  # was 10
$ws.Range('F5').Value = 1324  # was 1314
$ws.Range('F6').Value = 7649  # was 7643
$ws.Range('F8').Value = 109  # was 110
$ws.Range('F10').Value = 8446  # was 8440
$ws.Range('F14').Value = 5655  # was 5652
$ws.Range('F15').Value = 59  # was 58
$ws.Range('F16').Value = 2616  # was 2608
$ws.Range('F17').Value = 1136  # was 1131
$ws.Range('F18').Value = 4596  # was 4594
$ws.Range('B20').Value = '''2024-03-27'
$ws.Range('B20').Style = "Normal"
$ws.Range('C20').Value = '北京·跨越二次元ACG神级动漫世界巡回演唱会'
$ws.Range('D20').Value = '东三环北路36号 朝阳剧场'
$ws.Range('E20').Value = '2024.03.27 19:30-03.27 21:10'
$ws.Range('F20').Value = 104  # was 91
$ws.Range('G20').Value = 60  # was 238
$ws.Range('H20').Value = 'https://show.bilibili.com/platform/detail.html?id=81614'
$ws.Range('I20').Value = '//i0.hdslb.com/bfs/openplatform/202402/rAr8lSIU1706772309212.jpeg'
$ws.Range('B21').Value = '''2024-03-29'
$ws.Range('B21').Style = "Normal"
$ws.Range('C21').Value = '北京·2024图书市集春季场'
$ws.Range('D21').Value = '建国路郎家园6号 郎园Vintage'
$ws.Range('E21').Value = '2024.03.29 14:00-03.31 20:00'
$ws.Range('F21').Value = 33  # was 103
$ws.Range('G21').Value = 35.1  # was 60
$ws.Range('H21').Value = 'https://show.bilibili.com/platform/detail.html?id=81984'
$ws.Range('I21').Value = '//i1.hdslb.com/bfs/openplatform/202402/Zi09QvTC1708571966640.jpeg'
$ws.Range('B22').Value = '''2024-03-30'
$ws.Range('B22').Style = "Normal"
$ws.Range('C22').Value = '北京·「京都动画X春日计划」漫展感管弦室内乐音乐会 '
$ws.Range('D22').Value = '安慧里三区10号(北辰购物中心对面) 北京剧院'
$ws.Range('E22').Value = '2024.03.30 14:30-03.30 16:00'
$ws.Range('F22').Value = 123  # was 32
$ws.Range('G22').Value = 88  # was 35.1
$ws.Range('H22').Value = 'https://show.bilibili.com/platform/detail.html?id=82037'
$ws.Range('I22').Value = '//i2.hdslb.com/bfs/openplatform/202402/WSkVWKYr1708676572045.png'
$ws.Range('C23').Value = '北京·梦游园3.0代号鸢周年庆Only'
$ws.Range('D23').Value = '北花园路1号 超级蜂巢'
$ws.Range('E23').Value = '2024.03.30 10:00-03.30 17:00'
$ws.Range('F23').Value = 531  # was 121
$ws.Range('H23').Value = 'https://show.bilibili.com/platform/detail.html?id=81584'
$ws.Range('I23').Value = '//i2.hdslb.com/bfs/openplatform/202402/ASPwEB9W1706844758149.jpeg'
$ws.Range('C24').Value = '北京·集结 - 超级世代！ACGN LIVE 音乐节 4期'
$ws.Range('D24').Value = '日坛北路17号日坛公园北门对面 METAL BOX'
$ws.Range('E24').Value = '2024.03.30 18:00-03.30 22:00'
$ws.Range('F24').Value = 7  # was 528
$ws.Range('G24').Value = 108  # was 88
$ws.Range('H24').Value = 'https://show.bilibili.com/platform/detail.html?id=82457'
$ws.Range('I24').Value = '//i1.hdslb.com/bfs/openplatform/202403/3HZKJJSS1709778407525.jpeg'
$ws.Range('B25').Value = '''2024-04-04'
$ws.Range('B25').Style = "Normal"
$ws.Range('C25').Value = '北京·IDOx梦次元动漫游戏嘉年华3rd'
$ws.Range('D25').Value = '北京展览馆 北京展览馆'
$ws.Range('E25').Value = '2024.04.04 09:30-04.05 17:00'
$ws.Range('F25').Value = 3507  # was 6
$ws.Range('G25').Value = 80  # was 108
$ws.Range('H25').Value = 'https://show.bilibili.com/platform/detail.html?id=80825'
$ws.Range('I25').Value = '//i1.hdslb.com/bfs/openplatform/202402/P1YCG3MT1708329896103.jpeg'
$ws.Range('C26').Value = '北京·IDOx梦次元动漫游戏嘉年华3rd·配音演员 吴晛 专场活动'
$ws.Range('E26').Value = '2024.04.04 10:30-04.04 13:30'
$ws.Range('F26').Value = 55  # was 3487
$ws.Range('G26').Value = 198  # was 80
$ws.Range('H26').Value = 'https://show.bilibili.com/platform/detail.html?id=82490'
$ws.Range('I26').Value = '//i2.hdslb.com/bfs/openplatform/202403/ZVLr6IVF1709795299722.png'
$ws.Range('F28').Value = 21  # was 20
$ws.Range('F29').Value = 3000  # was 2976
$ws.Range('F30').Value = 343  # was 340
$ws.Range('F32').Value = 310  # was 307
$ws.Range('F34').Value = 351  # was 304
$ws.Range('F37').Value = 882  # was 879
$ws.Range('F39').Value = 1762  # was 1721
$ws.Range('F42').Value = 18  # was 16
$ws.Range('F43').Value = 2845  # was 2751
$ws.Range('F46').Value = 2284  # was 2283
$ws.Range('F47').Value = 30  # was 29
